$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 95
$prevRow = 94

# Copy the formatting (number format / style) from the row above so the
# new date cell picks up the same "yyyy-mm-dd hh:mm:ss" style already in
# use for column A, without registering a brand-new cell style.
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

# Plain numeric columns.
$ws.Cells.Item($newRow, 1).Value = 45460.2916666667
$ws.Cells.Item($newRow, 2).Value = 3300
$ws.Cells.Item($newRow, 3).Value = 0.709999978542328
$ws.Cells.Item($newRow, 4).Value = 0.709999978542328
$ws.Cells.Item($newRow, 5).Value = 0.709999978542328
$ws.Cells.Item($newRow, 6).Value = 0.709999978542328

# Column G in this sheet stores the "close" figure as text (shared
# string) rather than a number, matching every other row. Writing the
# value straight in would get auto-coerced to a number, so build it via
# a TEXT() formula and then flatten the formula down to its literal
# value - this keeps the result a genuine string cell without adding
# any new cell style.
$ws.Cells.Item($newRow, 7).Formula = "=TEXT(0.709999978542328,""0.000000000000000"")"
$ws.Cells.Item($newRow, 7).Copy()
$ws.Cells.Item($newRow, 7).PasteSpecial(-4163)

# Ticker column is already plain text.
$ws.Cells.Item($newRow, 8).Value = "BWZ.MI"

$excel.CutCopyMode = 0
